$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.364.06"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.846.50"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.41"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6308"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2955"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07706"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "1.867.87"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6828"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001004"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.88"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.125"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.407.30"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.39"
$ws.Range("E19").Value = "  -2.94%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.544"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.951"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.06"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1394"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.353"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05676"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.249"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.118"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7120"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "1.258.90"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01815"
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.779"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9121"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.217"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.20"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.14"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.061"
$ws.Range("E46").Value = "  -4.47%  "
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4039"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.069"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.681"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1123"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05734"
$ws.Range("E51").Value = "  -0.39%  "
